$d = $word.ActiveDocument
$p = $d.Paragraphs.Last
$r = $p.Range
$r.Collapse(0)
$xml = @'
<w:p><w:pPr><w:pStyle w:val="berschrift1"/></w:pPr><w:r><w:t>Small Turret</w:t></w:r></w:p><w:p><w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="0442E902" wp14:editId="671610A5"><wp:extent cx="1860870" cy="1648178"/><wp:effectExtent l="57150" t="0" r="0" b="0"/><wp:docPr id="16" name="Teilkreis 16"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="1860870" cy="1648178"/></a:xfrm><a:prstGeom prst="pie"><a:avLst><a:gd name="adj1" fmla="val 8564625"/><a:gd name="adj2" fmla="val 12950483"/></a:avLst></a:prstGeom><a:noFill/><a:ln w="127000"><a:solidFill><a:schemeClr val="tx1"/></a:solidFill><a:round/></a:ln></wps:spPr><wps:style><a:lnRef idx="2"><a:schemeClr val="accent1"><a:shade val="50000"/></a:schemeClr></a:lnRef><a:fillRef idx="1"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="0"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:inline></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shape w14:anchorId="4ECBB050" id="Teilkreis 16" o:spid="_x0000_s1026" style="width:146.55pt;height:129.8pt;visibility:visible;mso-wrap-style:square;mso-left-percent:-10001;mso-top-percent:-10001;mso-position-horizontal:absolute;mso-position-horizontal-relative:char;mso-position-vertical:absolute;mso-position-vertical-relative:line;mso-left-percent:-10001;mso-top-percent:-10001;v-text-anchor:middle" coordsize="1860870,1648178" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQDa2qCRsAIAAMkFAAAOAAAAZHJzL2Uyb0RvYy54bWysVEtv2zAMvg/YfxB0X/1YHm5QpwhadBhQ&#10;tMXaoWdVlmoPsqhJymu/fpSsONlW7DDsIosi+ZH8TPLictcrshHWdaBrWpzllAjNoen0a02/Pt18&#10;qChxnumGKdCipnvh6OXy/buLrVmIElpQjbAEQbRbbE1NW+/NIsscb0XP3BkYoVEpwfbMo2hfs8ay&#10;LaL3KivzfJZtwTbGAhfO4ev1oKTLiC+l4P5eSic8UTXF3Hw8bTxfwpktL9ji1TLTdjylwf4hi551&#10;GoOOUNfMM7K23R9QfcctOJD+jEOfgZQdF7EGrKbIf6vmsWVGxFqQHGdGmtz/g+V3m0fzYJGGrXEL&#10;h9dQxU7aPnwxP7KLZO1HssTOE46PRTXLqzlyylFXzCZVMa8CndnR3VjnPwnoSbjU1HShGLZgm1vn&#10;I1kN0azHrmDNt4IS2SvkfsMUqaazyaycpr9zYlWeWhXl+TSfVB9T1ASL8Q9xQwwNN51S8TcrTbaY&#10;ajnP8zwm4kB1TVAHw9hy4kpZghnU1O+KhPuLlYW1boYqlcZij6zFm98rEcCU/iIk6RrkqRxChYY+&#10;ojPOhfbFoGpZI4agU8ws9iQWMeYTKY2AAVliuiN2Angbe8gy2QdXEedhdE4c/M159IiRQfvRue80&#10;2LcqU1hVijzYH0gaqAksvUCzf7DEwjCNzvCbDhvkljn/wCy2ADYVrhR/j4dUgP8M0o2SFuyPt96D&#10;PU4FainZ4jjX1H1fMysoUZ81zst5MZmE+Y/CZDovUbCnmpdTjV73V4BNgE2J2cVrsPfqcJUW+mfc&#10;PKsQFVVMc4xdU+7tQbjyw5rB3cXFahXNcOYN87f60fAAHlgNvfq0e2bWpCnxOGB3cBj9NC0Do0fb&#10;4KlhtfYgOx+UR16TgPsiNk7abWEhncrR6riBlz8BAAD//wMAUEsDBBQABgAIAAAAIQBPoeiY3AAA&#10;AAUBAAAPAAAAZHJzL2Rvd25yZXYueG1sTI/BTsMwEETvSPyDtUjcqNMgQhviVFXVXEBCkPYDNvGS&#10;pI3tKHbb9O9ZuMBlNatZzbzNVpPpxZlG3zmrYD6LQJCtne5so2C/Kx4WIHxAq7F3lhRcycMqv73J&#10;MNXuYj/pXIZGcIj1KSpoQxhSKX3dkkE/cwNZ9r7caDDwOjZSj3jhcNPLOIoSabCz3NDiQJuW6mN5&#10;MgqS9+Jt97otqw+Hurge/D5+XmyVur+b1i8gAk3h7xh+8Bkdcmaq3MlqL3oF/Ej4nezFy8c5iIrF&#10;0zIBmWfyP33+DQAA//8DAFBLAQItABQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAAAAAAAAAAAAA&#10;AAAAAABbQ29udGVudF9UeXBlc10ueG1sUEsBAi0AFAAGAAgAAAAhADj9If/WAAAAlAEAAAsAAAAA&#10;AAAAAAAAAAAALwEAAF9yZWxzLy5yZWxzUEsBAi0AFAAGAAgAAAAhANraoJGwAgAAyQUAAA4AAAAA&#10;AAAAAAAAAAAALgIAAGRycy9lMm9Eb2MueG1sUEsBAi0AFAAGAAgAAAAhAE+h6JjcAAAABQEAAA8A&#10;AAAAAAAAAAAAAAAACgUAAGRycy9kb3ducmV2LnhtbFBLBQYAAAAABAAEAPMAAAATBgAAAAA=&#10;" path="m224554,1360975c-69027,1058176,-75418,612564,209383,303258l930435,824089,224554,1360975xe" filled="f" strokecolor="black [3213]" strokeweight="10pt"><v:path arrowok="t" o:connecttype="custom" o:connectlocs="224554,1360975;209383,303258;930435,824089;224554,1360975" o:connectangles="0,0,0,0"/><w10:anchorlock/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p><w:p><w:pPr><w:pStyle w:val="berschrift1"/></w:pPr><w:r><w:t>Small Turret Projectile</w:t></w:r></w:p><w:p><w:r><w:rPr><w:noProof/><w:highlight w:val="black"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="483EBD24" wp14:editId="07D2817A"><wp:extent cx="1152172" cy="992829"/><wp:effectExtent l="60642" t="53658" r="70803" b="70802"/><wp:docPr id="19" name="Gleichschenkliges Dreieck 19"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm rot="5400000"><a:off x="0" y="0"/><a:ext cx="1152172" cy="992829"/></a:xfrm><a:prstGeom prst="triangle"><a:avLst/></a:prstGeom><a:solidFill><a:schemeClr val="bg1"/></a:solidFill><a:ln w="127000"><a:solidFill><a:schemeClr val="bg1"/></a:solidFill><a:round/></a:ln></wps:spPr><wps:style><a:lnRef idx="2"><a:schemeClr val="accent1"><a:shade val="50000"/></a:schemeClr></a:lnRef><a:fillRef idx="1"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="0"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:inline></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shape w14:anchorId="4DE34952" id="Gleichschenkliges Dreieck 19" o:spid="_x0000_s1026" type="#_x0000_t5" style="width:90.7pt;height:78.2pt;rotation:90;visibility:visible;mso-wrap-style:square;mso-left-percent:-10001;mso-top-percent:-10001;mso-position-horizontal:absolute;mso-position-horizontal-relative:char;mso-position-vertical:absolute;mso-position-vertical-relative:line;mso-left-percent:-10001;mso-top-percent:-10001;v-text-anchor:middle" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQDtFd7XkQIAAK0FAAAOAAAAZHJzL2Uyb0RvYy54bWysVFFr2zAQfh/sPwi9r45NsjahTgktHYPS&#10;lqWjz4osxwJZp0lKnOzX7yTZbtMWBmN+EJLu7tN9n+/u8urQKrIX1knQJc3PJpQIzaGSelvSn0+3&#10;Xy4ocZ7piinQoqRH4ejV8vOny84sRAENqEpYgiDaLTpT0sZ7s8gyxxvRMncGRmg01mBb5vFot1ll&#10;WYforcqKyeRr1oGtjAUunMPbm2Sky4hf14L7h7p2whNVUszNx9XGdRPWbHnJFlvLTCN5nwb7hyxa&#10;JjU+OkLdMM/Izsp3UK3kFhzU/oxDm0FdSy4iB2STT96wWTfMiMgFxXFmlMn9P1h+v1+bR4sydMYt&#10;HG4Di0NtW2IB1ZpNJ+GL3DBbcojSHUfpxMETjpd5Pivy84ISjrb5vLgo5kHbLGEFTGOd/yagJWFT&#10;Um8l01sV6LEF2985n9wHt3DtQMnqVioVD6EkxLWyZM/wZ262ef/AiZfSpMNsivMh6RNrLKu/YljY&#10;6SplozRyeFEm7vxRiZCQ0j9ETWSF7IvI4g0641xonydTwyqREp9FPRP8GBGVioABuUbKI3YPcMp+&#10;wE4wvX8IFbHmx+D048ZnUganwWNEfBm0H4NbqcF+xEwhq/7l5D+IlKQJKm2gOj7aVEPYd87wW4n/&#10;/Y45/8gsthhe4tjwD7jUCvCfQb+jpAH7+6P74I+Vj1ZKOmzZkrpfO2YFJeq7xp6Y59Np6PF4mM7O&#10;CzzY15bNa4vetdeAhZTH7OI2+Hs1bGsL7TNOl1V4FU1Mc3y7pNzb4XDt0yjB+cTFahXdsK8N83d6&#10;bXgAD6qGmn46PDNrhuLHtrmHob3f1X/yDZEaVjsPtYzN8aJrrzfOhFg4/fwKQ+f1OXq9TNnlHwAA&#10;AP//AwBQSwMEFAAGAAgAAAAhALj79CDbAAAABQEAAA8AAABkcnMvZG93bnJldi54bWxMj81uwjAQ&#10;hO+VeAdrkXorDigglMZBiAiph15Kf86beInTxusoNiF9+5pe2stqVrOa+TbfTbYTIw2+daxguUhA&#10;ENdOt9woeHs9PmxB+ICssXNMCr7Jw66Y3eWYaXflFxpPoRExhH2GCkwIfSalrw1Z9AvXE0fv7AaL&#10;Ia5DI/WA1xhuO7lKko202HJsMNjTwVD9dbpYBavPevx4ei/L/pweq3V4NklaGqXu59P+EUSgKfwd&#10;ww0/okMRmSp3Ye1FpyA+En7nzVtvUhBVFNtlCrLI5X/64gcAAP//AwBQSwECLQAUAAYACAAAACEA&#10;toM4kv4AAADhAQAAEwAAAAAAAAAAAAAAAAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQA&#10;BgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAAAAAAAAAAAAAAC8BAABfcmVscy8ucmVsc1BLAQItABQA&#10;BgAIAAAAIQDtFd7XkQIAAK0FAAAOAAAAAAAAAAAAAAAAAC4CAABkcnMvZTJvRG9jLnhtbFBLAQIt&#10;ABQABgAIAAAAIQC4+/Qg2wAAAAUBAAAPAAAAAAAAAAAAAAAAAOsEAABkcnMvZG93bnJldi54bWxQ&#10;SwUGAAAAAAQABADzAAAA8wUAAAAA&#10;" fillcolor="white [3212]" strokecolor="white [3212]" strokeweight="10pt"><v:stroke joinstyle="round"/><w10:anchorlock/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r></w:p><w:p/>
'@
$r.InsertXML($xml)
